# Entrega RF0003 y RF0004
# Adds two new fields (COORDENADA_ESTE, COORDENADA_NORTE) to the "Datos"
# sheet of the infrastructure report template, inserted as new columns
# F and G (pushing the existing MAT_CONS_*/TIENE_CARTEL_*/EQUI_ACCESORIOS_*
# columns two places to the right, from F:U to H:W).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Datos")

# Insert two blank columns before the old column F. Excel shifts the
# existing cell data, column widths and the "plain" data validation
# (K2:M1048576 -> M2:O1048576) to the right automatically.
$ws.Columns("F:G").Insert()

# Fill in the headers for the two new fields.
$ws.Range("F1").Value = "COORDENADA_ESTE"
$ws.Range("G1").Value = "COORDENADA_NORTE"

# Match the column widths used for the new fields.
$ws.Columns("F").ColumnWidth = 19.71
$ws.Columns("G").ColumnWidth = 22

# Reflect the author's new active cell/selection on the sheet.
$ws.Range("H2").Select()
